# Update NATMI LR-pair data for Ccl4-Ccr5: the "Sending cluster" changes
# from ECs to MuSCs (new TPM-based values), and the MuSCs->MuSCs
# self-signalling row is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MuSCs -> ECs -------------------------------------------------
$ws.Cells.Item(2, 1).Value2 = "MuSCs"          # A2 Sending cluster
$ws.Cells.Item(2, 4).Value2 = "ECs"            # D2 Target cluster

$ws.Cells.Item(2, 7).Value2  = 0.2195956666666667   # G2
$ws.Cells.Item(2, 8).Value2  = 0.658787             # H2
$ws.Cells.Item(2, 13).Value2 = 0.001937666666666667 # M2
$ws.Cells.Item(2, 14).Value2 = 0.005813             # N2
$ws.Cells.Item(2, 15).Value2 = 0.0230007399171451   # O2
$ws.Cells.Item(2, 16).Value2 = 0.02300073991714511  # P2
$ws.Cells.Item(2, 17).Value2 = 0.0004255032034444444# Q2
$ws.Cells.Item(2, 18).Value2 = 0.003829528831       # R2
$ws.Cells.Item(2, 19).Value2 = 0.0230007399171451   # S2
$ws.Cells.Item(2, 20).Value2 = 0.02300073991714511  # T2

# --- Row 3: MuSCs -> FAPs -------------------------------------------------
$ws.Cells.Item(3, 1).Value2 = "MuSCs"          # A3 Sending cluster
$ws.Cells.Item(3, 4).Value2 = "FAPs"           # D3 Target cluster (unchanged text, same index shift)

$ws.Cells.Item(3, 7).Value2  = 0.2195956666666667   # G3
$ws.Cells.Item(3, 8).Value2  = 0.658787             # H3
$ws.Cells.Item(3, 13).Value2 = 0.082306             # M3
$ws.Cells.Item(3, 15).Value2 = 0.9769992600828549   # O3
$ws.Cells.Item(3, 16).Value2 = 0.976999260082855    # P3
$ws.Cells.Item(3, 17).Value2 = 0.01807404094066667  # Q3
$ws.Cells.Item(3, 18).Value2 = 0.162666368466       # R3
$ws.Cells.Item(3, 19).Value2 = 0.9769992600828549   # S3
$ws.Cells.Item(3, 20).Value2 = 0.976999260082855    # T3

# --- Row 4: MuSCs -> MuSCs self-pair, removed from the export -------------
$ws.Rows.Item(4).Delete()
